# Adding referral profile and modifying current profiles
#
# - Metadata sheet: bump the generated "Date" value and update the
#   concept "Count" from 3 to 1 (only the OPD row remains on Concepts).
# - Concepts sheet: remove the IPD and EMG rows, keeping only the OPD row.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-08-01T12:25:19+00:00"

# "Count" is stored as text in this metadata table (like every other
# property on this sheet), so force a text number format before writing
# the digit - otherwise Excel would happily re-interpret "1" as a number.
$wsMeta.Range("B22").NumberFormat = "@"
$wsMeta.Range("B22").Value = "1"

$wsConcepts = $wb.Worksheets.Item("Concepts")
$wsConcepts.Rows.Item(3).Delete()
$wsConcepts.Rows.Item(3).Delete()
